$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.712.87'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.303.52'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'108.98"
$ws.Range('E5').Value = '  +12.40%  '
$ws.Range('D6').Value = "'270.66"
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').Value = "'0.618"
$ws.Range('E7').Value = '  -1.77%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = "'0.618"
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').Value = "'48.04"
$ws.Range('E10').Value = '  +6.87%  '
$ws.Range('D11').Value = "'0.0939"
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = "'8.67"
$ws.Range('E12').Value = '  +9.21%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = "'15.72"
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('D15').Value = '2.645.14'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = "'0.857"
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').Value = '2.297.90'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').Value = '43.688.04'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('D21').Value = "'72.26"
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').Value = "'2.52"
$ws.Range('E22').Value = '  +8.24%  '
$ws.Range('D23').Value = "'233.79"
$ws.Range('E23').Value = '  -2.66%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = "'9.53"
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = "'2.88"
$ws.Range('E25').Value = '  +13.47%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'42.97"
$ws.Range('E28').Value = '  +10.70%  '
$ws.Range('D29').Value = "'3.46"
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = "'177.35"
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = "'21.77"
$ws.Range('E32').Value = '  -2.46%  '
$ws.Range('D33').Value = "'0.0921"
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('D34').Value = "'5.66"
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = "'4.81"
$ws.Range('E35').Value = '  +7.21%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = "'0.127"
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = "'0.112"
$ws.Range('E37').Value = '  +1.72%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = "'4.03"
$ws.Range('E38').Value = '  +20.01%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.0356"
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').Value = "'0.237"
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = "'70.01"
$ws.Range('E42').Value = '  +11.77%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = "'1.38"
$ws.Range('E44').Value = '  -3.13%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').Value = "'12.16"
$ws.Range('E45').Value = '  -1.24%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = "'5.43"
$ws.Range('E46').Value = '  +1.41%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = "'8.74"
$ws.Range('E47').Value = '  -4.40%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = "'0.102"
$ws.Range('E48').Value = '  -1.85%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = "'99.62"
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = "'1.22"
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = "'0.452"
$ws.Range('E51').Value = '  +7.25%  '
